$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-12 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-13 Saturday", 2)

$d.Content.Find.Execute("534÷5=106, 4", $true, $false, $false, $false, $false, $true, 1, $false, "378÷3=126, 0", 2)
$d.Content.Find.Execute("919÷7=131, 2", $true, $false, $false, $false, $false, $true, 1, $false, "791÷4=197, 3", 2)
$d.Content.Find.Execute("646÷4=161, 2", $true, $false, $false, $false, $false, $true, 1, $false, "731÷6=121, 5", 2)
$d.Content.Find.Execute("350÷5=70, 0", $true, $false, $false, $false, $false, $true, 1, $false, "899÷5=179, 4", 2)
$d.Content.Find.Execute("227÷2=113, 1", $true, $false, $false, $false, $false, $true, 1, $false, "348÷3=116, 0", 2)

$d.Content.Find.Execute("407÷8=50, 7", $true, $false, $false, $false, $false, $true, 1, $false, "584÷4=146, 0", 2)
$d.Content.Find.Execute("273÷9=30, 3", $true, $false, $false, $false, $false, $true, 1, $false, "617÷9=68, 5", 2)
$d.Content.Find.Execute("386÷6=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "960÷3=320, 0", 2)
$d.Content.Find.Execute("262÷7=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "389÷7=55, 4", 2)
$d.Content.Find.Execute("245÷6=40, 5", $true, $false, $false, $false, $false, $true, 1, $false, "904÷6=150, 4", 2)

$d.Content.Find.Execute("921÷8=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "469÷5=93, 4", 2)
$d.Content.Find.Execute("320÷8=40, 0", $true, $false, $false, $false, $false, $true, 1, $false, "751÷6=125, 1", 2)
$d.Content.Find.Execute("664÷7=94, 6", $true, $false, $false, $false, $false, $true, 1, $false, "853÷5=170, 3", 2)
$d.Content.Find.Execute("211÷4=52, 3", $true, $false, $false, $false, $false, $true, 1, $false, "296÷7=42, 2", 2)
$d.Content.Find.Execute("423÷4=105, 3", $true, $false, $false, $false, $false, $true, 1, $false, "133÷2=66, 1", 2)

$d.Content.Find.Execute("928÷8=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "626÷5=125, 1", 2)
$d.Content.Find.Execute("188÷8=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "469÷2=234, 1", 2)
$d.Content.Find.Execute("180÷2=90, 0", $true, $false, $false, $false, $false, $true, 1, $false, "204÷6=34, 0", 2)
$d.Content.Find.Execute("497÷6=82, 5", $true, $false, $false, $false, $false, $true, 1, $false, "586÷8=73, 2", 2)
$d.Content.Find.Execute("479÷8=59, 7", $true, $false, $false, $false, $false, $true, 1, $false, "974÷2=487, 0", 2)

$d.Content.Find.Execute("444÷7=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "757÷4=189, 1", 2)
$d.Content.Find.Execute("230÷4=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "117÷9=13, 0", 2)
$d.Content.Find.Execute("336÷6=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "607÷8=75, 7", 2)
$d.Content.Find.Execute("883÷6=147, 1", $true, $false, $false, $false, $false, $true, 1, $false, "360÷5=72, 0", 2)
$d.Content.Find.Execute("679÷8=84, 7", $true, $false, $false, $false, $false, $true, 1, $false, "923÷8=115, 3", 2)
